# Logboek.xlsx update:
# - Add a new log entry on row 7 (date 22/03/2022, work description, time)
# - Update the active selection to B7

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Row 7: date, time (shared string 11), work description (shared string 12)
# Copy the date-style formatting from the row above so A7 reuses the
# existing date style instead of Excel minting a brand new number format.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A7").Value = (Get-Date -Year 2022 -Month 3 -Day 22 -Hour 0 -Minute 0 -Second 0)

$ws.Range("C7").Value = "6u30"

$ws.Range("B7").Value = "De timer werkt perfect samen met de andere code, (ik heb timer interrupt hiervoor moeten gebruiken). En ik heb ook de start en reset knop hun gewenste functie kunnen geven (ook via interrupts). Het versturen van het reset signaal is ook klaar. "

# Row 7 wraps onto two lines, like the other multi-line entries above it
$ws.Rows.Item(7).RowHeight = 28.8

# Update the active selection to B7
$ws.Range("B7").Select()
